# Add two new user rows (Mac-Addresses) to the master-user_detail_h sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows below the last data row (30) so they inherit the
# same per-column formatting (email column fill flag, is_active left align)
# that every existing data row already carries.
$ws.Rows("31:32").Insert()

# Row 32: John Doe (entered first, so its strings land earlier in the shared table)
$ws.Range("A32").Value = 110031
$ws.Range("B32").Value = 9317596767
$ws.Range("C32").Value = "John Doe"
$ws.Range("D32").Value = "john.doe@xyz.com"
$ws.Range("E32").Value = 818876431
$ws.Range("F32").Value = "ACT"
$ws.Range("G32").Value = "eng"
$ws.Range("H32").Value = "PWD"
$ws.Range("I32").Value = $true
$ws.Range("J32").Value = "superadmin"
$ws.Range("K32").Value = "now()"
$ws.Range("L32").Value = "now()"

# Row 31: Jane Smith
$ws.Range("A31").Value = 110030
$ws.Range("B31").Value = 9317596768
$ws.Range("C31").Value = "Jane Smith"
$ws.Range("D31").Value = "jane.smith@xyz.com"
$ws.Range("E31").Value = 818876432
$ws.Range("F31").Value = "ACT"
$ws.Range("G31").Value = "eng"
$ws.Range("H31").Value = "PWD"
$ws.Range("I31").Value = $true
$ws.Range("J31").Value = "superadmin"
$ws.Range("K31").Value = "now()"
$ws.Range("L31").Value = "now()"

# Selection moved to F30 in the final file
$ws.Range("F30").Select()
